$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new team-record columns, reusing the header style
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill the team record (1996 Mets: 71-91-0) for every player row (2-42)
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 29).Value = 71
    $ws.Cells.Item($r, 30).Value = 91
    $ws.Cells.Item($r, 31).Value = 0
}
